$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 2 correction: "unnamed: 1_level_1" -> "total" ---
$ws.Cells.Item(2, 2).Value = "total"

# --- Data block: rows 4-38 shift up (two section-header rows removed,
#     trailing rows 39-40 dropped) ---
$ws.Cells.Item(4, 1).Value = "               brasil"
$ws.Cells.Item(4, 2).Value = 0.6906713719105717
$ws.Cells.Item(4, 3).Value = 0.6783216753484271
$ws.Cells.Item(4, 4).Value = 0.7332794004464949
$ws.Cells.Item(4, 5).Value = 0.7417404274906203
$ws.Cells.Item(4, 6).Value = 0.9130288365589294
$ws.Cells.Item(4, 7).Value = 1.354377001987363
$ws.Cells.Item(4, 8).Value = 0.9567611433089162

$ws.Cells.Item(5, 1).Value = "urbana"
$ws.Cells.Item(5, 2).Value = 0.7006971388827247
$ws.Cells.Item(5, 3).Value = 0.699592194727477
$ws.Cells.Item(5, 4).Value = 0.7440967620751812
$ws.Cells.Item(5, 5).Value = 0.7804105407800235
$ws.Cells.Item(5, 6).Value = 0.9078572465655519
$ws.Cells.Item(5, 7).Value = 1.403045057963147
$ws.Cells.Item(5, 8).Value = 0.9556825476124285

$ws.Cells.Item(6, 1).Value = "rural"
$ws.Cells.Item(6, 2).Value = 1.493477155056306
$ws.Cells.Item(6, 3).Value = 1.430355237574423
$ws.Cells.Item(6, 4).Value = 1.613697579993041
$ws.Cells.Item(6, 5).Value = 1.481761452261324
$ws.Cells.Item(6, 6).Value = 1.812018080366819
$ws.Cells.Item(6, 7).Value = 3.288647252027231
$ws.Cells.Item(6, 8).Value = 1.824880185534455

$ws.Cells.Item(7, 1).Value = "norte"
$ws.Cells.Item(7, 2).Value = 3.397100745969429
$ws.Cells.Item(7, 3).Value = 3.365967680970556
$ws.Cells.Item(7, 4).Value = 3.478946217880297
$ws.Cells.Item(7, 5).Value = 3.32607778025787
$ws.Cells.Item(7, 6).Value = 3.539612428526985
$ws.Cells.Item(7, 7).Value = 5.338905859801104
$ws.Cells.Item(7, 8).Value = 3.561735960963752

$ws.Cells.Item(8, 1).Value = "rondônia"
$ws.Cells.Item(8, 2).Value = 2.799733739902242
$ws.Cells.Item(8, 3).Value = 2.868446494254792
$ws.Cells.Item(8, 4).Value = 2.961043930156892
$ws.Cells.Item(8, 5).Value = 3.123268870469983
$ws.Cells.Item(8, 6).Value = 3.077066425510282
$ws.Cells.Item(8, 7).Value = 5.157025240230544
$ws.Cells.Item(8, 8).Value = 3.291830750909619

$ws.Cells.Item(9, 1).Value = "acre"
$ws.Cells.Item(9, 2).Value = 11.45605862855805
$ws.Cells.Item(9, 3).Value = 11.00045835451473
$ws.Cells.Item(9, 4).Value = 12.04514383327654
$ws.Cells.Item(9, 5).Value = 15.01738722462592
$ws.Cells.Item(9, 6).Value = 10.62666032589985
$ws.Cells.Item(9, 7).Value = 21.78260503750971
$ws.Cells.Item(9, 8).Value = 10.35426930007584

$ws.Cells.Item(10, 1).Value = "amazonas"
$ws.Cells.Item(10, 2).Value = 4.276446206546409
$ws.Cells.Item(10, 3).Value = 4.109031838505013
$ws.Cells.Item(10, 4).Value = 4.559244039751499
$ws.Cells.Item(10, 5).Value = 4.969378041082455
$ws.Cells.Item(10, 6).Value = 4.696878190377769
$ws.Cells.Item(10, 7).Value = 13.02961323369042
$ws.Cells.Item(10, 8).Value = 4.660098951634417

$ws.Cells.Item(11, 1).Value = "roraima"
$ws.Cells.Item(11, 2).Value = 4.376312156228948
$ws.Cells.Item(11, 3).Value = 4.602195572510944
$ws.Cells.Item(11, 4).Value = 4.571757029331318
$ws.Cells.Item(11, 5).Value = 6.689363317662469
$ws.Cells.Item(11, 6).Value = 4.122989748565124
$ws.Cells.Item(11, 7).Value = 5.79703144382406
$ws.Cells.Item(11, 8).Value = 4.240160147517347

$ws.Cells.Item(12, 1).Value = "pará"
$ws.Cells.Item(12, 2).Value = 8.557475742034605
$ws.Cells.Item(12, 3).Value = 8.423666259714535
$ws.Cells.Item(12, 4).Value = 8.796665365895445
$ws.Cells.Item(12, 5).Value = 8.80645373067177
$ws.Cells.Item(12, 6).Value = 8.528865016720108
$ws.Cells.Item(12, 7).Value = 11.50948961168409
$ws.Cells.Item(12, 8).Value = 8.53007616913708

$ws.Cells.Item(13, 1).Value = "amapá"
$ws.Cells.Item(13, 2).Value = 8.76438314984467
$ws.Cells.Item(13, 3).Value = 8.442079567605811
$ws.Cells.Item(13, 4).Value = 9.557949019076688
$ws.Cells.Item(13, 5).Value = 11.71006372948207
$ws.Cells.Item(13, 6).Value = 7.78310021773393
$ws.Cells.Item(13, 7).Value = 13.77523728489339
$ws.Cells.Item(13, 8).Value = 7.612810033239758

$ws.Cells.Item(14, 1).Value = "tocantins"
$ws.Cells.Item(14, 2).Value = 3.128747707857803
$ws.Cells.Item(14, 3).Value = 2.496381216179003
$ws.Cells.Item(14, 4).Value = 3.929552657737186
$ws.Cells.Item(14, 5).Value = 3.488123437081244
$ws.Cells.Item(14, 6).Value = 3.478831430140461
$ws.Cells.Item(14, 7).Value = 6.010486912922588
$ws.Cells.Item(14, 8).Value = 3.347795780931667

$ws.Cells.Item(15, 1).Value = "nordeste"
$ws.Cells.Item(15, 2).Value = 1.671660636326493
$ws.Cells.Item(15, 3).Value = 1.593930277547595
$ws.Cells.Item(15, 4).Value = 1.783072970454697
$ws.Cells.Item(15, 5).Value = 1.899642536410373
$ws.Cells.Item(15, 6).Value = 1.821026294444619
$ws.Cells.Item(15, 7).Value = 2.960686417109789
$ws.Cells.Item(15, 8).Value = 1.857054248575162

$ws.Cells.Item(16, 1).Value = "maranhão"
$ws.Cells.Item(16, 2).Value = 8.004064990087455
$ws.Cells.Item(16, 3).Value = 7.455689188990165
$ws.Cells.Item(16, 4).Value = 8.625977385613508
$ws.Cells.Item(16, 5).Value = 8.044260003219371
$ws.Cells.Item(16, 6).Value = 8.7956806147176
$ws.Cells.Item(16, 7).Value = 8.152895174527059
$ws.Cells.Item(16, 8).Value = 9.29283692321688

$ws.Cells.Item(17, 1).Value = "piauí"
$ws.Cells.Item(17, 2).Value = 5.660274782198973
$ws.Cells.Item(17, 3).Value = 5.048366860295426
$ws.Cells.Item(17, 4).Value = 6.269864590741218
$ws.Cells.Item(17, 5).Value = 7.36594166844625
$ws.Cells.Item(17, 6).Value = 5.712985905990586
$ws.Cells.Item(17, 7).Value = 14.69479919235982
$ws.Cells.Item(17, 8).Value = 5.410540234974242

$ws.Cells.Item(18, 1).Value = "ceará"
$ws.Cells.Item(18, 2).Value = 4.07463890643841
$ws.Cells.Item(18, 3).Value = 3.776373361904417
$ws.Cells.Item(18, 4).Value = 4.453970591167084
$ws.Cells.Item(18, 5).Value = 3.787229230185474
$ws.Cells.Item(18, 6).Value = 4.517394012113559
$ws.Cells.Item(18, 7).Value = 8.773344958130345
$ws.Cells.Item(18, 8).Value = 4.454158692189392

$ws.Cells.Item(19, 1).Value = "rio grande do norte"
$ws.Cells.Item(19, 2).Value = 9.048653178928317
$ws.Cells.Item(19, 3).Value = 8.92830790118979
$ws.Cells.Item(19, 4).Value = 9.312784715293839
$ws.Cells.Item(19, 5).Value = 9.088833841871866
$ws.Cells.Item(19, 6).Value = 9.259214182847009
$ws.Cells.Item(19, 7).Value = 11.49274301576103
$ws.Cells.Item(19, 8).Value = 9.479094241294092

$ws.Cells.Item(20, 1).Value = "paraíba"
$ws.Cells.Item(20, 2).Value = 5.584940896054857
$ws.Cells.Item(20, 3).Value = 5.297019968193077
$ws.Cells.Item(20, 4).Value = 5.901246380905792
$ws.Cells.Item(20, 5).Value = 7.141147034409292
$ws.Cells.Item(20, 6).Value = 5.313242978114399
$ws.Cells.Item(20, 7).Value = 6.921975428972462
$ws.Cells.Item(20, 8).Value = 5.519254357478987

$ws.Cells.Item(21, 1).Value = "pernambuco"
$ws.Cells.Item(21, 2).Value = 2.756252961321327
$ws.Cells.Item(21, 3).Value = 2.995372884932341
$ws.Cells.Item(21, 4).Value = 2.645941607340438
$ws.Cells.Item(21, 5).Value = 3.110004816489339
$ws.Cells.Item(21, 6).Value = 2.952977817077343
$ws.Cells.Item(21, 7).Value = 6.480505524780043
$ws.Cells.Item(21, 8).Value = 2.992352985805404

$ws.Cells.Item(22, 1).Value = "alagoas"
$ws.Cells.Item(22, 2).Value = 8.569540049541555
$ws.Cells.Item(22, 3).Value = 7.86108646280681
$ws.Cells.Item(22, 4).Value = 9.279432653716409
$ws.Cells.Item(22, 5).Value = 12.31669023711709
$ws.Cells.Item(22, 6).Value = 8.367293657515834
$ws.Cells.Item(22, 7).Value = 11.53595941246714
$ws.Cells.Item(22, 8).Value = 8.799296430414188

$ws.Cells.Item(23, 1).Value = "sergipe"
$ws.Cells.Item(23, 2).Value = 5.551153259693107
$ws.Cells.Item(23, 3).Value = 6.400862244665305
$ws.Cells.Item(23, 4).Value = 5.147064379353981
$ws.Cells.Item(23, 5).Value = 7.104140619491832
$ws.Cells.Item(23, 6).Value = 5.453065152438271
$ws.Cells.Item(23, 7).Value = 11.49576082867021
$ws.Cells.Item(23, 8).Value = 5.652140203977453

$ws.Cells.Item(24, 1).Value = "bahia"
$ws.Cells.Item(24, 2).Value = 2.717221721983847
$ws.Cells.Item(24, 3).Value = 2.579792554493431
$ws.Cells.Item(24, 4).Value = 2.952622441996626
$ws.Cells.Item(24, 5).Value = 3.595076394770903
$ws.Cells.Item(24, 6).Value = 3.029431688076834
$ws.Cells.Item(24, 7).Value = 4.382751739252297
$ws.Cells.Item(24, 8).Value = 2.988580144804423

$ws.Cells.Item(25, 1).Value = "sudeste"
$ws.Cells.Item(25, 2).Value = 0.941662504558406
$ws.Cells.Item(25, 3).Value = 0.9457984392997021
$ws.Cells.Item(25, 4).Value = 1.00976525304489
$ws.Cells.Item(25, 5).Value = 1.099213896262621
$ws.Cells.Item(25, 6).Value = 1.155153958562273
$ws.Cells.Item(25, 7).Value = 1.914819447415016
$ws.Cells.Item(25, 8).Value = 1.212881961510793

$ws.Cells.Item(26, 1).Value = "minas gerais"
$ws.Cells.Item(26, 2).Value = 1.872325235677548
$ws.Cells.Item(26, 3).Value = 1.849936163565545
$ws.Cells.Item(26, 4).Value = 1.976440487362646
$ws.Cells.Item(26, 5).Value = 2.254525246015203
$ws.Cells.Item(26, 6).Value = 2.016901690751915
$ws.Cells.Item(26, 7).Value = 3.058427934336676
$ws.Cells.Item(26, 8).Value = 2.059223880320989

$ws.Cells.Item(27, 1).Value = "espírito santo"
$ws.Cells.Item(27, 2).Value = 3.820866609919103
$ws.Cells.Item(27, 3).Value = 3.712779475081575
$ws.Cells.Item(27, 4).Value = 4.341633040877873
$ws.Cells.Item(27, 5).Value = 4.591037829601102
$ws.Cells.Item(27, 6).Value = 4.278704966822766
$ws.Cells.Item(27, 7).Value = 7.558234882328136
$ws.Cells.Item(27, 8).Value = 4.604065764650269

$ws.Cells.Item(28, 1).Value = "rio de janeiro"
$ws.Cells.Item(28, 2).Value = 2.27049841959801
$ws.Cells.Item(28, 3).Value = 2.352072327325425
$ws.Cells.Item(28, 4).Value = 2.357219853733232
$ws.Cells.Item(28, 5).Value = 2.88915175267143
$ws.Cells.Item(28, 6).Value = 2.698260701037416
$ws.Cells.Item(28, 7).Value = 3.934321143841288
$ws.Cells.Item(28, 8).Value = 2.851093449126763

$ws.Cells.Item(29, 1).Value = "são paulo"
$ws.Cells.Item(29, 2).Value = 1.299701116255888
$ws.Cells.Item(29, 3).Value = 1.301741418166295
$ws.Cells.Item(29, 4).Value = 1.421451595419515
$ws.Cells.Item(29, 5).Value = 1.445960525245666
$ws.Cells.Item(29, 6).Value = 1.783191834517231
$ws.Cells.Item(29, 7).Value = 3.258243534346394
$ws.Cells.Item(29, 8).Value = 1.91819370117347

$ws.Cells.Item(30, 1).Value = "sul"
$ws.Cells.Item(30, 2).Value = 1.370830733887644
$ws.Cells.Item(30, 3).Value = 1.348079384956707
$ws.Cells.Item(30, 4).Value = 1.449532280402686
$ws.Cells.Item(30, 5).Value = 1.490045449338303
$ws.Cells.Item(30, 6).Value = 1.781783107137463
$ws.Cells.Item(30, 7).Value = 2.784242319045342
$ws.Cells.Item(30, 8).Value = 1.956187363983048

$ws.Cells.Item(31, 1).Value = "paraná"
$ws.Cells.Item(31, 2).Value = 2.196571755331435
$ws.Cells.Item(31, 3).Value = 2.169872153179015
$ws.Cells.Item(31, 4).Value = 2.388615996717964
$ws.Cells.Item(31, 5).Value = 2.302139229197608
$ws.Cells.Item(31, 6).Value = 2.982054944115066
$ws.Cells.Item(31, 7).Value = 4.599176270582457
$ws.Cells.Item(31, 8).Value = 3.189947169764456

$ws.Cells.Item(32, 1).Value = "santa catarina"
$ws.Cells.Item(32, 2).Value = 3.68014068075559
$ws.Cells.Item(32, 3).Value = 3.454352007377039
$ws.Cells.Item(32, 4).Value = 3.966095533273659
$ws.Cells.Item(32, 5).Value = 3.914371636584687
$ws.Cells.Item(32, 6).Value = 3.561015319746415
$ws.Cells.Item(32, 7).Value = 5.901140428594061
$ws.Cells.Item(32, 8).Value = 3.947379896049716

$ws.Cells.Item(33, 1).Value = "rio grande do sul"
$ws.Cells.Item(33, 2).Value = 1.640927281200707
$ws.Cells.Item(33, 3).Value = 1.715846034705353
$ws.Cells.Item(33, 4).Value = 1.6442889951643
$ws.Cells.Item(33, 5).Value = 1.675858580962141
$ws.Cells.Item(33, 6).Value = 2.5914885417739
$ws.Cells.Item(33, 7).Value = 4.149754775020637
$ws.Cells.Item(33, 8).Value = 2.789020564164013

$ws.Cells.Item(34, 1).Value = "centro-oeste"
$ws.Cells.Item(34, 2).Value = 1.411729464675901
$ws.Cells.Item(34, 3).Value = 1.459668995638385
$ws.Cells.Item(34, 4).Value = 1.536584000617008
$ws.Cells.Item(34, 5).Value = 1.600351256671607
$ws.Cells.Item(34, 6).Value = 1.723639418381639
$ws.Cells.Item(34, 7).Value = 3.060821377747543
$ws.Cells.Item(34, 8).Value = 1.760305977710202

$ws.Cells.Item(35, 1).Value = "mato grosso do sul"
$ws.Cells.Item(35, 2).Value = 3.209940761680517
$ws.Cells.Item(35, 3).Value = 3.305770252653994
$ws.Cells.Item(35, 4).Value = 3.314514397448112
$ws.Cells.Item(35, 5).Value = 3.165965067798285
$ws.Cells.Item(35, 6).Value = 3.731660377210302
$ws.Cells.Item(35, 7).Value = 5.834325126106514
$ws.Cells.Item(35, 8).Value = 3.749490030997287

$ws.Cells.Item(36, 1).Value = "mato grosso"
$ws.Cells.Item(36, 2).Value = 2.77654037066981
$ws.Cells.Item(36, 3).Value = 3.231149079747463
$ws.Cells.Item(36, 4).Value = 2.647635607631215
$ws.Cells.Item(36, 5).Value = 4.089308202968292
$ws.Cells.Item(36, 6).Value = 3.616419374638118
$ws.Cells.Item(36, 7).Value = 5.94953139165065
$ws.Cells.Item(36, 8).Value = 3.612095045633005

$ws.Cells.Item(37, 1).Value = "goiás"
$ws.Cells.Item(37, 2).Value = 2.212800153550007
$ws.Cells.Item(37, 3).Value = 2.163605390757528
$ws.Cells.Item(37, 4).Value = 2.574080208747711
$ws.Cells.Item(37, 5).Value = 2.526912274958109
$ws.Cells.Item(37, 6).Value = 2.582384948681153
$ws.Cells.Item(37, 7).Value = 4.659424951902824
$ws.Cells.Item(37, 8).Value = 2.692751340005902

$ws.Cells.Item(38, 1).Value = "distrito federal"
$ws.Cells.Item(38, 2).Value = 3.731641278420229
$ws.Cells.Item(38, 3).Value = 3.675302146628663
$ws.Cells.Item(38, 4).Value = 4.007827837110157
$ws.Cells.Item(38, 5).Value = 3.611562673602033
$ws.Cells.Item(38, 6).Value = 4.325510589784596
$ws.Cells.Item(38, 7).Value = 8.143091484409547
$ws.Cells.Item(38, 8).Value = 4.293185569947054

# --- Remove now-unused trailing rows 39 and 40 ---
$ws.Range("A39:H40").Delete()

$ws.Range("A1").Select()